$wb = $excel.ActiveWorkbook

# ---- Templates ----
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "Templates"

$ws.Range("A1").Value = "title"
$ws.Range("B1").Value = "description"
$ws.Range("C1").Value = "version"
$ws.Range("D1").Value = "tags"
$ws.Range("A1:D1").Interior.Color = 65535

$ws.Range("A2").Value = "Template1"
$ws.Range("B2").Value = "Template1 description"
$ws.Range("C2").Value = "v1"
$ws.Range("D2").Value = "Template 1 tags"

$ws.Range("A3").Value = "Template2"
$ws.Range("B3").Value = "Template2 description"
$ws.Range("C3").Value = "v2"
$ws.Range("D3").Value = "Template 2 tags"

$ws.Rows("1:1").Select()

# ---- MultiDocs ----
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "MultiDocs"

$ws.Range("A1").Value = "folder"
$ws.Range("B1").Value = "contact"
$ws.Range("C1").Value = "company"
$ws.Range("D1").Value = "prospect"
$ws.Range("E1").Value = "task"
$ws.Range("F1").Value = "case"
$ws.Range("G1").Value = "tags"
$ws.Range("A1:G1").Interior.Color = 65535

$ws.Range("A2").Value = "Root directory"
$ws.Range("B2").Value = "contact1"
$ws.Range("C2").Value = "company1"
$ws.Range("D2").Value = "prospect1"
$ws.Range("E2").Value = "task1"
$ws.Range("F2").Value = "case1"
$ws.Range("G2").Value = "tag1"

$ws.Range("A3").Value = "Word Templates"
$ws.Range("B3").Value = "contact2"
$ws.Range("C3").Value = "company2"
$ws.Range("D3").Value = "prospect2"
$ws.Range("E3").Value = "task2"
$ws.Range("F3").Value = "case2"
$ws.Range("G3").Value = "tag2"

$ws.Rows("1:1").Select()

# ---- FeedbackForms ----
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "FeedbackForms"

$ws.Range("A1").Value = "title"
$ws.Range("B1").Value = "pages"
$ws.Range("C1").Value = "reportEmail"
$ws.Range("D1").Value = "Description"
$ws.Range("E1").Value = "welcomeMessage"
$ws.Range("F1").Value = "confirmationMessage"
$ws.Range("A1:F1").Interior.Color = 65535

$ws.Range("A2").Value = "Form title 1"
$ws.Range("B2").Value = "'10"
$ws.Range("C2").Value = "abc@gmail.com"
$ws.Range("D2").Value = "Form 1 description"
$ws.Range("E2").Value = "Form 1 welcome message"
$ws.Range("F2").Value = "Form1 confirmation message"

$ws.Range("A3").Value = "Form title 2"
$ws.Range("B3").Value = "'15"
$ws.Range("C3").Value = "def@gmail.com"
$ws.Range("D3").Value = "Form 2 description"
$ws.Range("E3").Value = "Form 2 welcome message"
$ws.Range("F3").Value = "Form2 confirmation message"

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:abc@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:def@gmail.com")

$ws.Range("F3").Select()

# ---- Tasks ----
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "Tasks"

$ws.Range("A1").Value = "title"
$ws.Range("B1").Value = "autoextend"
$ws.Range("C1").Value = "status"
$ws.Range("D1").Value = "completionpercentage"
$ws.Range("E1").Value = "type"
$ws.Range("F1").Value = "priority"
$ws.Range("G1").Value = "deal"
$ws.Range("H1").Value = "case"
$ws.Range("I1").Value = "tags"
$ws.Range("J1").Value = "description"
$ws.Range("K1").Value = "contact"
$ws.Range("L1").Value = "company"
$ws.Range("M1").Value = "identifier"
$ws.Range("A1:M1").Interior.Color = 65535

$ws.Range("A2").Value = "Task1"
$ws.Range("B2").Value = "Extend deadline by 1 day"
$ws.Range("C2").Value = "Open"
$ws.Range("D2").Value = "'50"
$ws.Range("E2").Value = "Call"
$ws.Range("F2").Value = "High"
$ws.Range("G2").Value = "deal1"
$ws.Range("H2").Value = "case1"
$ws.Range("I2").Value = "tag1"
$ws.Range("J2").Value = "tagdesc1"
$ws.Range("K2").Value = "contact1"
$ws.Range("L2").Value = "company1"
$ws.Range("M2").Value = "identifier1"

$ws.Range("A3").Value = "Task2"
$ws.Range("B3").Value = "Extend deadline by 30 days"
$ws.Range("C3").Value = "Complete"
$ws.Range("D3").Value = "'75"
$ws.Range("E3").Value = "Training"
$ws.Range("F3").Value = "Normal"
$ws.Range("G3").Value = "deal2"
$ws.Range("H3").Value = "case2"
$ws.Range("I3").Value = "tag2"
$ws.Range("J3").Value = "tagdesc2"
$ws.Range("K3").Value = "contact2"
$ws.Range("L3").Value = "company2"
$ws.Range("M3").Value = "identifier2"

$ws.Range("I19").Select()
